# Slovenia Prva Liga data refresh
# - Insert one new match as the new row 128 (pushes the two previously-last rows down by one)
# - Append six brand-new matches as rows 131-136
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 128 used to be the second-to-last row. Insert a fresh row above it so the two
# existing rows (126, 127 match-ids) shift down to rows 129-130, matching the diff.
$ws.Rows.Item(128).Insert()

# Helper to (re)apply the two custom cell styles used throughout this sheet:
#  - column A: bold, centered, thin border (style index 1 in the original workbook)
#  - column E: custom date/time display format (style index 2 in the original workbook)
function Set-IdStyle($addr) {
  $r = $ws.Range($addr)
  $r.Font.Bold = $true
  $r.HorizontalAlignment = -4108
  $r.VerticalAlignment = -4160
  $r.Borders.LineStyle = 1
}
function Set-DateStyle($addr) {
  $r = $ws.Range($addr)
  $r.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

foreach ($rowNum in 128..136) {
  Set-IdStyle("A$rowNum")
  Set-DateStyle("E$rowNum")
}

# Row 128
$ws.Range("A128").Value = 126
$ws.Range("B128").Value = 6816447
$ws.Range("C128").Value = "Slovenia Prva Liga"
$ws.Range("D128").Value = "Slovenia Prva Liga"
$ws.Range("E128").Value = 45368.375
$ws.Range("F128").Value = "NK Rogaska"
$ws.Range("G128").Value = "NS Mura"
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = "D"
$ws.Range("K128").Value = 2.4
$ws.Range("L128").Value = 3.2
$ws.Range("M128").Value = 2.625
$ws.Range("N128").Value = 2.1
$ws.Range("O128").Value = 3.2
$ws.Range("P128").Value = 3.25
$ws.Range("Q128").Value = -0.25
$ws.Range("R128").Value = 1.9
$ws.Range("S128").Value = 1.9
$ws.Range("T128").Value = 2.25
$ws.Range("U128").Value = 1.825
$ws.Range("V128").Value = 1.975
$ws.Range("W128").Value = -1
$ws.Range("X128").Value = 2.2
$ws.Range("Y128").Value = -1
$ws.Range("Z128").Value = -0.5
$ws.Range("AA128").Value = 0.45
$ws.Range("AB128").Value = -1
$ws.Range("AC128").Value = 0.9750000000000001

# Row 129
$ws.Range("A129").Value = 127
$ws.Range("B129").Value = 6814429
$ws.Range("C129").Value = "Slovenia Prva Liga"
$ws.Range("D129").Value = "Slovenia Prva Liga"
$ws.Range("E129").Value = 45368.45833333334
$ws.Range("F129").Value = "NK Celje"
$ws.Range("G129").Value = "NK Radomlje"
$ws.Range("H129").Value = 2
$ws.Range("I129").Value = 1
$ws.Range("J129").Value = "H"
$ws.Range("K129").Value = 1.285
$ws.Range("L129").Value = 5
$ws.Range("M129").Value = 7.5
$ws.Range("N129").Value = 1.333
$ws.Range("O129").Value = 4.5
$ws.Range("P129").Value = 7
$ws.Range("Q129").Value = -1.5
$ws.Range("R129").Value = 1.975
$ws.Range("S129").Value = 1.825
$ws.Range("T129").Value = 2.75
$ws.Range("U129").Value = 1.95
$ws.Range("V129").Value = 1.85
$ws.Range("W129").Value = 0.333
$ws.Range("X129").Value = -1
$ws.Range("Y129").Value = -1
$ws.Range("Z129").Value = -1
$ws.Range("AA129").Value = 0.825
$ws.Range("AB129").Value = 0.475
$ws.Range("AC129").Value = -0.5

# Row 130
$ws.Range("A130").Value = 128
$ws.Range("B130").Value = 6814430
$ws.Range("C130").Value = "Slovenia Prva Liga"
$ws.Range("D130").Value = "Slovenia Prva Liga"
$ws.Range("E130").Value = 45368.5625
$ws.Range("F130").Value = "NK Maribor"
$ws.Range("G130").Value = "FC Koper"
$ws.Range("H130").Value = 3
$ws.Range("I130").Value = 1
$ws.Range("J130").Value = "H"
$ws.Range("K130").Value = 1.666
$ws.Range("L130").Value = 3.6
$ws.Range("M130").Value = 4.2
$ws.Range("N130").Value = 1.727
$ws.Range("O130").Value = 3.5
$ws.Range("P130").Value = 4
$ws.Range("Q130").Value = -0.75
$ws.Range("R130").Value = 1.975
$ws.Range("S130").Value = 1.825
$ws.Range("T130").Value = 2.5
$ws.Range("U130").Value = 1.85
$ws.Range("V130").Value = 1.95
$ws.Range("W130").Value = 0.7270000000000001
$ws.Range("X130").Value = -1
$ws.Range("Y130").Value = -1
$ws.Range("Z130").Value = 0.9750000000000001
$ws.Range("AA130").Value = -1
$ws.Range("AB130").Value = 0.8500000000000001
$ws.Range("AC130").Value = -1

# Row 131
$ws.Range("A131").Value = 129
$ws.Range("B131").Value = 7907436
$ws.Range("C131").Value = "Slovenia Prva Liga"
$ws.Range("D131").Value = "Slovenia Prva Liga"
$ws.Range("E131").Value = 45374.67708333334
$ws.Range("F131").Value = "FC Koper"
$ws.Range("G131").Value = "NS Mura"
$ws.Range("H131").Value = 2
$ws.Range("I131").Value = 1
$ws.Range("J131").Value = "H"
$ws.Range("K131").Value = 1.85
$ws.Range("L131").Value = 3.4
$ws.Range("M131").Value = 3.8
$ws.Range("N131").Value = 1.95
$ws.Range("O131").Value = 3.25
$ws.Range("P131").Value = 3.5
$ws.Range("Q131").Value = -0.5
$ws.Range("R131").Value = 2
$ws.Range("S131").Value = 1.8
$ws.Range("T131").Value = 2.5
$ws.Range("U131").Value = 1.975
$ws.Range("V131").Value = 1.825
$ws.Range("W131").Value = 0.95
$ws.Range("X131").Value = -1
$ws.Range("Y131").Value = -1
$ws.Range("Z131").Value = 1
$ws.Range("AA131").Value = -1
$ws.Range("AB131").Value = 0.9750000000000001
$ws.Range("AC131").Value = -1

# Row 132
$ws.Range("A132").Value = 130
$ws.Range("B132").Value = 6816446
$ws.Range("C132").Value = "Slovenia Prva Liga"
$ws.Range("D132").Value = "Slovenia Prva Liga"
$ws.Range("E132").Value = 45380.45833333334
$ws.Range("F132").Value = "NK Rogaska"
$ws.Range("G132").Value = "NK Bravo"
$ws.Range("K132").Value = 2.8
$ws.Range("L132").Value = 3.5
$ws.Range("M132").Value = 2.1
$ws.Range("N132").Value = 2.75
$ws.Range("O132").Value = 3.5
$ws.Range("P132").Value = 2.1
$ws.Range("Q132").Value = 0.25
$ws.Range("R132").Value = 1.85
$ws.Range("S132").Value = 1.95
$ws.Range("T132").Value = 2.5
$ws.Range("U132").Value = 1.95
$ws.Range("V132").Value = 1.85
$ws.Range("W132").Value = 0
$ws.Range("X132").Value = 0
$ws.Range("Y132").Value = 0
$ws.Range("Z132").Value = 0
$ws.Range("AA132").Value = 0

# Row 133
$ws.Range("A133").Value = 131
$ws.Range("B133").Value = 7977922
$ws.Range("C133").Value = "Slovenia Prva Liga"
$ws.Range("D133").Value = "Slovenia Prva Liga"
$ws.Range("E133").Value = 45380.5625
$ws.Range("F133").Value = "NK Maribor"
$ws.Range("G133").Value = "NK Radomlje"
$ws.Range("K133").Value = 1.285
$ws.Range("L133").Value = 5.5
$ws.Range("M133").Value = 6.5
$ws.Range("N133").Value = 1.4
$ws.Range("O133").Value = 5
$ws.Range("P133").Value = 5
$ws.Range("Q133").Value = -1.25
$ws.Range("R133").Value = 1.85
$ws.Range("S133").Value = 1.95
$ws.Range("T133").Value = 2.75
$ws.Range("U133").Value = 1.85
$ws.Range("V133").Value = 1.95
$ws.Range("W133").Value = 0
$ws.Range("X133").Value = 0
$ws.Range("Y133").Value = 0
$ws.Range("Z133").Value = 0
$ws.Range("AA133").Value = 0

# Row 134
$ws.Range("A134").Value = 132
$ws.Range("B134").Value = 7977924
$ws.Range("C134").Value = "Slovenia Prva Liga"
$ws.Range("D134").Value = "Slovenia Prva Liga"
$ws.Range("E134").Value = 45380.67708333334
$ws.Range("F134").Value = "NK Domzale"
$ws.Range("G134").Value = "NS Mura"
$ws.Range("K134").Value = 2
$ws.Range("L134").Value = 3.4
$ws.Range("M134").Value = 3.1
$ws.Range("N134").Value = 2.3
$ws.Range("O134").Value = 3.4
$ws.Range("P134").Value = 2.6
$ws.Range("Q134").Value = 0
$ws.Range("R134").Value = 1.75
$ws.Range("S134").Value = 2.05
$ws.Range("T134").Value = 2.5
$ws.Range("U134").Value = 1.9
$ws.Range("V134").Value = 1.9
$ws.Range("W134").Value = 0
$ws.Range("X134").Value = 0
$ws.Range("Y134").Value = 0
$ws.Range("Z134").Value = 0
$ws.Range("AA134").Value = 0

# Row 135
$ws.Range("A135").Value = 133
$ws.Range("B135").Value = 7977921
$ws.Range("C135").Value = "Slovenia Prva Liga"
$ws.Range("D135").Value = "Slovenia Prva Liga"
$ws.Range("E135").Value = 45381.45833333334
$ws.Range("F135").Value = "Olimpija Ljubljana"
$ws.Range("G135").Value = "FC Koper"
$ws.Range("K135").Value = 1.5
$ws.Range("L135").Value = 3.6
$ws.Range("M135").Value = 6
$ws.Range("N135").Value = 1.5
$ws.Range("O135").Value = 3.6
$ws.Range("P135").Value = 5.75
$ws.Range("Q135").Value = -1
$ws.Range("R135").Value = 1.9
$ws.Range("S135").Value = 1.9
$ws.Range("T135").Value = 2.75
$ws.Range("U135").Value = 1.95
$ws.Range("V135").Value = 1.85
$ws.Range("W135").Value = 0
$ws.Range("X135").Value = 0
$ws.Range("Y135").Value = 0
$ws.Range("Z135").Value = 0
$ws.Range("AA135").Value = 0

# Row 136
$ws.Range("A136").Value = 134
$ws.Range("B136").Value = 7977923
$ws.Range("C136").Value = "Slovenia Prva Liga"
$ws.Range("D136").Value = "Slovenia Prva Liga"
$ws.Range("E136").Value = 45381.5625
$ws.Range("F136").Value = "NK Celje"
$ws.Range("G136").Value = "NK Aluminij"
$ws.Range("K136").Value = 1.2
$ws.Range("L136").Value = 6.5
$ws.Range("M136").Value = 8
$ws.Range("N136").Value = 1.2
$ws.Range("O136").Value = 6.5
$ws.Range("P136").Value = 8
$ws.Range("Q136").Value = -1.75
$ws.Range("R136").Value = 1.825
$ws.Range("S136").Value = 1.975
$ws.Range("T136").Value = 3
$ws.Range("U136").Value = 1.85
$ws.Range("V136").Value = 1.95
$ws.Range("W136").Value = 0
$ws.Range("X136").Value = 0
$ws.Range("Y136").Value = 0
$ws.Range("Z136").Value = 0
$ws.Range("AA136").Value = 0

Write-Output "Slovenia Prva Liga sheet updated: inserted row 128, appended rows 131-136."
